$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $seatCell = $ws.Cells.Item($r, 8)   # Column H - Seat Type
    $rankCell = $ws.Cells.Item($r, 9)   # Column I - Rank

    $seatVal = $seatCell.Text
    if ($seatVal -ne $null -and $seatVal -ne "") {
        $seatStr = $seatVal.ToString()
        $newSeatVal = $seatStr.ToUpper()
        if ($newSeatVal.EndsWith(":")) {
            $newSeatVal = $newSeatVal.Substring(0, $newSeatVal.Length - 1)
        }
        if (-not $newSeatVal.Equals($seatStr)) {
            $seatCell.Value = $newSeatVal
        }
    }

    $rankVal = $rankCell.Text
    if ($rankVal -ne $null -and $rankVal -ne "") {
        $rankStr = $rankVal.ToString()
        if ($rankStr.EndsWith(",")) {
            $newRankVal = $rankStr.Substring(0, $rankStr.Length - 1)
            $rankCell.Value = $newRankVal
        }
    }
}
